$wb = $excel.ActiveWorkbook

$resSheet = $wb.Worksheets.Item("RES installed")
$mainSheet = $wb.Worksheets.Item("Main")

$resSheet.Range("C2").Value = 8
$resSheet.Range("C3").Value = 10
$resSheet.Range("C4").Value = 6
$resSheet.Range("C5").Value = 6
$resSheet.Range("C6").Value = 3

$excel.CalculateFullRebuild()

$resSheet.Range("C3").Select()
$resSheet.Activate()

$mainSheet.Range("B4").Select()
